$wb = $excel.ActiveWorkbook

# Rename the "object" sheet to "car"
$carSheet = $wb.Worksheets.Item("object")
$carSheet.Name = "car"

# Strip the "car:" prefix from the header row of the car sheet
$carSheet.Range("B1").Value = "brand"
$carSheet.Range("C1").Value = "model"
$carSheet.Range("D1").Value = "fuel"
$carSheet.Range("E1").Value = "year"

# Re-apply the existing date format so both date cells collapse onto the
# same (already-registered) style record instead of keeping a stray
# duplicate around
$carSheet.Range("E2").NumberFormat = "yyyy\-mm\-dd"
$carSheet.Range("E3").NumberFormat = "yyyy\-mm\-dd"

# Strip the "computer:" prefix from the header row of the computers sheet
$computersSheet = $wb.Worksheets.Item("computers")
$computersSheet.Range("B1").Value = "cpu"
$computersSheet.Range("C1").Value = "gpu"
$computersSheet.Range("D1").Value = "storage"

# Match the recorded cursor/selection position on each sheet
$carSheet.Activate()
$carSheet.Range("E7").Select()

$computersSheet.Activate()
$computersSheet.Range("B1").Select()
